$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2 to 546) all currently hold the date serial value 45204
# (2023-10-05) and need to be updated to 45205 (2023-10-06).
$lastRow = 546

$range = $ws.Range("C2:C" + $lastRow)
$range.Value = 45205
